$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.964.69'
$ws.Range('E2').Value = '  +4.58%  '
$ws.Range('D3').Value = '2.679.88'
$ws.Range('E3').Value = '  +8.08%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '113.66'
$ws.Range('E5').Value = '  +9.62%  '
$ws.Range('D6').Value = '326.23'
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.553'
$ws.Range('E9').Value = '  +3.93%  '
$ws.Range('D10').Value = '41.00'
$ws.Range('E10').Value = '  +6.32%  '
$ws.Range('D11').Value = '20.11'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '0.0824'
$ws.Range('E12').Value = '  +3.50%  '
$ws.Range('D14').Value = '7.37'
$ws.Range('E14').Value = '  +5.17%  '
$ws.Range('D15').Value = '3.091.64'
$ws.Range('E15').Value = '  +7.82%  '
$ws.Range('D16').Value = '2.673.26'
$ws.Range('E16').Value = '  +7.63%  '
$ws.Range('D17').Value = '0.876'
$ws.Range('E17').Value = '  +6.60%  '
$ws.Range('D18').Value = '49.879.35'
$ws.Range('E18').Value = '  +4.56%  '
$ws.Range('D19').Value = '13.18'
$ws.Range('E19').Value = '  +4.23%  '
$ws.Range('D20').Value = '6.80'
$ws.Range('E20').Value = '  +4.60%  '
$ws.Range('D21').Value = '2.89'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '0.0₃0959'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('D23').Value = '278.29'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '71.83'
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('E25').Value = '  +3.66%  '
$ws.Range('D26').Value = '26.91'
$ws.Range('E26').Value = '  +5.10%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '10.12'
$ws.Range('E28').Value = '  +6.12%  '
$ws.Range('E29').Value = '  +1.21%  '
$ws.Range('D30').Value = '36.24'
$ws.Range('E30').Value = '  +5.73%  '
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('D32').Value = '50.34'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('E33').Value = '  +4.78%  '
$ws.Range('D34').Value = '19.55'
$ws.Range('E34').Value = '  +3.71%  '
$ws.Range('D35').Value = '0.0811'
$ws.Range('E35').Value = '  +6.06%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '5.06'
$ws.Range('E37').Value = '  +13.08%  '
$ws.Range('D38').Value = '2.08'
$ws.Range('E38').Value = '  +7.80%  '
$ws.Range('D39').Value = '3.16'
$ws.Range('E39').Value = '  +11.08%  '
$ws.Range('D40').Value = '124.68'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('D42').Value = '22.78'
$ws.Range('E42').Value = '  +5.76%  '
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('E44').Value = '  +6.76%  '
$ws.Range('D45').Value = '2.112.72'
$ws.Range('E45').Value = '  +6.27%  '
$ws.Range('D46').Value = '3.31'
$ws.Range('E46').Value = '  +6.13%  '
$ws.Range('E47').Value = '  +15.28%  '
$ws.Range('D48').Value = '2.06'
$ws.Range('E48').Value = '  +9.25%  '
$ws.Range('D49').Value = '9.05'
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('D50').Value = '5.37'
$ws.Range('E50').Value = '  +6.28%  '
$ws.Range('D51').Value = '59.74'
$ws.Range('E51').Value = '  +7.10%  '
